# Auto-applied numeric corrections to Leve profit-tracking sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H:N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 18183330
$ws.Cells.Item(100, 9).Value = 20001514
$ws.Cells.Item(100, 10).Value = 1500
$ws.Cells.Item(100, 11).Value = 20001514
$ws.Cells.Item(100, 12).Value = 1500
$ws.Cells.Item(100, 13).Value = -20000973
$ws.Cells.Item(100, 14).Value = -2582

$ws.Cells.Item(137, 8).Value = 1362632
$ws.Cells.Item(137, 10).Value = 3008.5
$ws.Cells.Item(137, 12).Value = 9025.5
$ws.Cells.Item(137, 14).Value = -14125.5

$ws.Cells.Item(138, 8).Value = 6691.87
$ws.Cells.Item(138, 9).Value = 811.5454999999999
$ws.Cells.Item(138, 10).Value = 8350.423000000001
$ws.Cells.Item(138, 11).Value = 2434.6365
$ws.Cells.Item(138, 12).Value = 25051.269
$ws.Cells.Item(138, 13).Value = 2705.3635
$ws.Cells.Item(138, 14).Value = -35331.269


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1390.1052
$ws.Cells.Item(61, 9).Value = 1324.2941
$ws.Cells.Item(61, 11).Value = 1324.2941
$ws.Cells.Item(61, 13).Value = -1112.2941

$ws.Cells.Item(74, 8).Value = 4425.433
$ws.Cells.Item(74, 9).Value = 5129.737
$ws.Cells.Item(74, 11).Value = 5129.737
$ws.Cells.Item(74, 13).Value = -4255.737

$ws.Cells.Item(77, 8).Value = 4425.433
$ws.Cells.Item(77, 9).Value = 5129.737
$ws.Cells.Item(77, 11).Value = 25648.685
$ws.Cells.Item(77, 13).Value = -21280.685

$ws.Cells.Item(88, 8).Value = 9526481
$ws.Cells.Item(88, 9).Value = 11113727
$ws.Cells.Item(88, 10).Value = 3000
$ws.Cells.Item(88, 11).Value = 11113727
$ws.Cells.Item(88, 12).Value = 3000
$ws.Cells.Item(88, 13).Value = -11113321
$ws.Cells.Item(88, 14).Value = -3812

$ws.Cells.Item(91, 8).Value = 9526481
$ws.Cells.Item(91, 9).Value = 11113727
$ws.Cells.Item(91, 10).Value = 3000
$ws.Cells.Item(91, 11).Value = 11113727
$ws.Cells.Item(91, 12).Value = 3000
$ws.Cells.Item(91, 13).Value = -11112323
$ws.Cells.Item(91, 14).Value = -5808

$ws.Cells.Item(102, 8).Value = 1495
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 14).Value = -5244

$ws.Cells.Item(125, 8).Value = 41805.625
$ws.Cells.Item(125, 10).Value = 41805.625
$ws.Cells.Item(125, 12).Value = 41805.625
$ws.Cells.Item(125, 14).Value = -51645.625

$ws.Cells.Item(132, 8).Value = 2281.524
$ws.Cells.Item(132, 9).Value = 1328.5
$ws.Cells.Item(132, 10).Value = 7999.6665
$ws.Cells.Item(132, 11).Value = 3985.5
$ws.Cells.Item(132, 12).Value = 23998.9995
$ws.Cells.Item(132, 13).Value = -1455.5
$ws.Cells.Item(132, 14).Value = -29058.9995

$ws.Cells.Item(136, 8).Value = 1390.1052
$ws.Cells.Item(136, 9).Value = 1324.2941
$ws.Cells.Item(136, 11).Value = 3972.8823
$ws.Cells.Item(136, 13).Value = -1422.8823


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2450.5833
$ws.Cells.Item(86, 9).Value = 2050.875
$ws.Cells.Item(86, 10).Value = 3250
$ws.Cells.Item(86, 11).Value = 2050.875
$ws.Cells.Item(86, 12).Value = 3250
$ws.Cells.Item(86, 13).Value = -927.875
$ws.Cells.Item(86, 14).Value = -5496

$ws.Cells.Item(89, 8).Value = 2450.5833
$ws.Cells.Item(89, 9).Value = 2050.875
$ws.Cells.Item(89, 10).Value = 3250
$ws.Cells.Item(89, 11).Value = 10254.375
$ws.Cells.Item(89, 12).Value = 16250
$ws.Cells.Item(89, 13).Value = -4638.375
$ws.Cells.Item(89, 14).Value = -27482

$ws.Cells.Item(133, 8).Value = 100000
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).ClearContents()

$ws.Cells.Item(134, 8).Value = 2727.6
$ws.Cells.Item(134, 9).Value = 1992.8334
$ws.Cells.Item(134, 10).Value = 5666.6665
$ws.Cells.Item(134, 11).Value = 5978.5002
$ws.Cells.Item(134, 12).Value = 16999.9995
$ws.Cells.Item(134, 13).Value = -3443.5002
$ws.Cells.Item(134, 14).Value = -22069.9995


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2626.2593
$ws.Cells.Item(31, 9).Value = 1186.1428
$ws.Cells.Item(31, 10).Value = 7666.6665
$ws.Cells.Item(31, 11).Value = 1186.1428
$ws.Cells.Item(31, 12).Value = 7666.6665
$ws.Cells.Item(31, 13).Value = -891.1428000000001
$ws.Cells.Item(31, 14).Value = -8256.666499999999

$ws.Cells.Item(34, 8).Value = 2626.2593
$ws.Cells.Item(34, 9).Value = 1186.1428
$ws.Cells.Item(34, 10).Value = 7666.6665
$ws.Cells.Item(34, 11).Value = 1186.1428
$ws.Cells.Item(34, 12).Value = 7666.6665
$ws.Cells.Item(34, 13).Value = -984.1428000000001
$ws.Cells.Item(34, 14).Value = -8070.6665

$ws.Cells.Item(58, 8).Value = 2600.8784
$ws.Cells.Item(58, 9).Value = 1575.6508
$ws.Cells.Item(58, 10).Value = 8472.637000000001
$ws.Cells.Item(58, 11).Value = 1575.6508
$ws.Cells.Item(58, 12).Value = 8472.637000000001
$ws.Cells.Item(58, 13).Value = -1372.6508
$ws.Cells.Item(58, 14).Value = -8878.637000000001

$ws.Cells.Item(105, 8).Value = 1359.826
$ws.Cells.Item(105, 9).Value = 1132.4
$ws.Cells.Item(105, 10).Value = 1786.25
$ws.Cells.Item(105, 11).Value = 1132.4
$ws.Cells.Item(105, 12).Value = 1786.25
$ws.Cells.Item(105, 13).Value = 614.5999999999999
$ws.Cells.Item(105, 14).Value = -5280.25

$ws.Cells.Item(132, 8).Value = 2395.9062
$ws.Cells.Item(132, 9).Value = 1380.7273
$ws.Cells.Item(132, 11).Value = 4142.1819
$ws.Cells.Item(132, 13).Value = -1612.1819

$ws.Cells.Item(134, 8).Value = 4927.5483
$ws.Cells.Item(134, 9).Value = 6291.579
$ws.Cells.Item(134, 10).Value = 2767.8333
$ws.Cells.Item(134, 11).Value = 18874.737
$ws.Cells.Item(134, 12).Value = 8303.499899999999
$ws.Cells.Item(134, 13).Value = -16339.737
$ws.Cells.Item(134, 14).Value = -13373.4999

$ws.Cells.Item(136, 8).Value = 2600.8784
$ws.Cells.Item(136, 9).Value = 1575.6508
$ws.Cells.Item(136, 10).Value = 8472.637000000001
$ws.Cells.Item(136, 11).Value = 4726.9524
$ws.Cells.Item(136, 12).Value = 25417.911
$ws.Cells.Item(136, 13).Value = -2176.9524
$ws.Cells.Item(136, 14).Value = -30517.911


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 8334607
$ws.Cells.Item(113, 9).Value = 1218.5555
$ws.Cells.Item(113, 10).Value = 20834690
$ws.Cells.Item(113, 11).Value = 3655.6665
$ws.Cells.Item(113, 12).Value = 62504070
$ws.Cells.Item(113, 13).Value = -1485.6665
$ws.Cells.Item(113, 14).Value = -62508410


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5849.838
$ws.Cells.Item(70, 9).Value = 5557.25
$ws.Cells.Item(70, 10).Value = 6760.1113
$ws.Cells.Item(70, 11).Value = 5557.25
$ws.Cells.Item(70, 12).Value = 6760.1113
$ws.Cells.Item(70, 13).Value = -5287.25
$ws.Cells.Item(70, 14).Value = -7300.1113

$ws.Cells.Item(73, 8).Value = 5849.838
$ws.Cells.Item(73, 9).Value = 5557.25
$ws.Cells.Item(73, 10).Value = 6760.1113
$ws.Cells.Item(73, 11).Value = 5557.25
$ws.Cells.Item(73, 12).Value = 6760.1113
$ws.Cells.Item(73, 13).Value = -4621.25
$ws.Cells.Item(73, 14).Value = -8632.1113

$ws.Cells.Item(80, 8).Value = 50002384
$ws.Cells.Item(80, 9).Value = 62502228
$ws.Cells.Item(80, 11).Value = 62502228
$ws.Cells.Item(80, 13).Value = -62501230

$ws.Cells.Item(83, 8).Value = 50002384
$ws.Cells.Item(83, 9).Value = 62502228
$ws.Cells.Item(83, 11).Value = 312511140
$ws.Cells.Item(83, 13).Value = -312506148

$ws.Cells.Item(132, 8).Value = 2767.3333
$ws.Cells.Item(132, 9).Value = 998.6667
$ws.Cells.Item(132, 10).Value = 6304.6665
$ws.Cells.Item(132, 11).Value = 2996.0001
$ws.Cells.Item(132, 12).Value = 18913.9995
$ws.Cells.Item(132, 13).Value = -466.0001000000002
$ws.Cells.Item(132, 14).Value = -23973.9995


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1576.7778
$ws.Cells.Item(16, 9).Value = 1398.875
$ws.Cells.Item(16, 10).Value = 3000
$ws.Cells.Item(16, 11).Value = 1398.875
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = -1228.875
$ws.Cells.Item(16, 14).Value = -3340

$ws.Cells.Item(42, 8).Value = 34631.4
$ws.Cells.Item(42, 10).Value = 34631.4
$ws.Cells.Item(42, 12).Value = 34631.4
$ws.Cells.Item(42, 14).Value = -35757.4

$ws.Cells.Item(46, 8).Value = 1938.5264
$ws.Cells.Item(46, 9).Value = 1608.3334
$ws.Cells.Item(46, 10).Value = 2504.5715
$ws.Cells.Item(46, 11).Value = 1608.3334
$ws.Cells.Item(46, 12).Value = 2504.5715
$ws.Cells.Item(46, 13).Value = -1420.3334
$ws.Cells.Item(46, 14).Value = -2880.5715

$ws.Cells.Item(49, 8).Value = 34631.4
$ws.Cells.Item(49, 10).Value = 34631.4
$ws.Cells.Item(49, 12).Value = 34631.4
$ws.Cells.Item(49, 14).Value = -34925.4

$ws.Cells.Item(122, 8).Value = 6596.1665
$ws.Cells.Item(122, 9).Value = 2626
$ws.Cells.Item(122, 11).Value = 7878
$ws.Cells.Item(122, 13).Value = -5428

$ws.Cells.Item(132, 8).Value = 3810.9143
$ws.Cells.Item(132, 9).Value = 1428.7826
$ws.Cells.Item(132, 10).Value = 8376.666999999999
$ws.Cells.Item(132, 11).Value = 4286.3478
$ws.Cells.Item(132, 12).Value = 25130.001
$ws.Cells.Item(132, 13).Value = -1756.3478
$ws.Cells.Item(132, 14).Value = -30190.001

$ws.Cells.Item(136, 8).Value = 3800.2188
$ws.Cells.Item(136, 9).Value = 1655.6923
$ws.Cells.Item(136, 10).Value = 5267.5264
$ws.Cells.Item(136, 11).Value = 4967.0769
$ws.Cells.Item(136, 12).Value = 15802.5792
$ws.Cells.Item(136, 13).Value = -2417.0769
$ws.Cells.Item(136, 14).Value = -20902.5792


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 19611012
$ws.Cells.Item(132, 9).Value = 2039.3334
$ws.Cells.Item(132, 10).Value = 41671104
$ws.Cells.Item(132, 11).Value = 6118.0002
$ws.Cells.Item(132, 12).Value = 125013312
$ws.Cells.Item(132, 13).Value = -3588.0002
$ws.Cells.Item(132, 14).Value = -125018372

$ws.Cells.Item(135, 8).Value = 33880.445
$ws.Cells.Item(135, 10).Value = 33880.445
$ws.Cells.Item(135, 12).Value = 33880.445
$ws.Cells.Item(135, 14).Value = -44020.445

$ws.Cells.Item(136, 8).Value = 4136.6553
$ws.Cells.Item(136, 9).Value = 2316.8096
$ws.Cells.Item(136, 10).Value = 8913.75
$ws.Cells.Item(136, 11).Value = 6950.4288
$ws.Cells.Item(136, 12).Value = 26741.25
$ws.Cells.Item(136, 13).Value = -4400.4288
$ws.Cells.Item(136, 14).Value = -31841.25
